$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row (row 1) from "_old"/"_new" suffixes to "_FV2404"/"_FV2410"
$ws.Range("A1").Value = "Segmentname_FV2404"
$ws.Range("B1").Value = "Segmentgruppe_FV2404"
$ws.Range("C1").Value = "Segment_FV2404"
$ws.Range("D1").Value = "Datenelement_FV2404"
$ws.Range("E1").Value = "Segment ID_FV2404"
$ws.Range("F1").Value = "Code_FV2404"
$ws.Range("G1").Value = "Qualifier_FV2404"
$ws.Range("H1").Value = "Beschreibung_FV2404"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2404"
$ws.Range("J1").Value = "Bedingung_FV2404"

$ws.Range("L1").Value = "Segmentname_FV2410"
$ws.Range("M1").Value = "Segmentgruppe_FV2410"
$ws.Range("N1").Value = "Segment_FV2410"
$ws.Range("O1").Value = "Datenelement_FV2410"
$ws.Range("P1").Value = "Segment ID_FV2410"
$ws.Range("Q1").Value = "Code_FV2410"
$ws.Range("R1").Value = "Qualifier_FV2410"
$ws.Range("S1").Value = "Beschreibung_FV2410"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2410"
$ws.Range("U1").Value = "Bedingung_FV2410"

# Temporarily stash the header row's existing formatting so that converting the
# range into a Table does not "capture" it as a brand-new header dxf (which a
# freshly inserted, still-unstyled header does not need).
$headerRange = $ws.Range("A1:U1")
$backupRange = $ws.Range("A1000:U1000")
$headerRange.Copy()
$backupRange.PasteSpecial(-4122)
$headerRange.ClearFormats()

# Turn the data range into an Excel Table ("Table1") spanning the full used range
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U90"), $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

# Restore the original header formatting and discard the backup
$backupRange.Copy()
$headerRange.PasteSpecial(-4122)
$backupRange.Clear()
$excel.CutCopyMode = $false

# Freeze the header row
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
